# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计" and fill
#    it with the per-fund holding detail for the new quarter.
# 2) Insert a new row at the top of the data in "总计" summarising the new
#    quarter (holdings count / market value), pushing the existing
#    "2021-Q4" row down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted right after "2021-Q4"
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$ws.Name = "2022-Q1"

# Header row (bold, bordered, centered - matches the other sheets' header style)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: col A is a numeric row index (0-based), styled like the header;
# columns D/E/F/G hold numeric-looking figures that must stay TEXT (as in
# the source data), column H is a genuine number.
$rows = @(
    @("003719", "易方达标普医疗保健指数(QDII-LOF) 美元", "0.51", "94.20", "1.51", "0.0077", 9),
    @("161126", "易方达标普医疗保健指数(QDII-LOF) 人民币", "0.51", "94.20", "1.51", "0.0077", 9),
    @("519981", "长信美国标准普尔100等权重指数增强(QDII)", "0.47", "84.16", "0.86", "0.0040", 10),
    @("011706", "长信美国标准普尔100等权重指数增强(QDII) - 美元", "0.47", "84.16", "0.86", "0.0040", 10)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = 2 + $r
    $data = $rows[$r]

    $idxCell = $ws.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $bCell = $ws.Cells.Item($rowNum, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $data[0]

    $ws.Cells.Item($rowNum, 3).Value = $data[1]

    $dCell = $ws.Cells.Item($rowNum, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $data[2]

    $eCell = $ws.Cells.Item($rowNum, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $data[3]

    $fCell = $ws.Cells.Item($rowNum, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $data[4]

    $gCell = $ws.Cells.Item($rowNum, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $data[5]

    $ws.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to "总计" (now pushed to the 3rd tab)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$aCell = $totalSheet.Cells.Item(2, 1)
$aCell.Value = 0
$aCell.Font.Bold = $true
$aCell.HorizontalAlignment = -4108
$aCell.VerticalAlignment = -4160
$aCell.Borders.LineStyle = 1

$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 4
$totalSheet.Cells.Item(2, 4).Value = 0.02

# The row that used to be row 2 ("2021-Q4") is now row 3; its running
# index (column A) needs to move from 0 to 1 to stay sequential.
$totalSheet.Cells.Item(3, 1).Value = 1
